$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17 (diff hunk @@ -1468,22 +1468,22 @@)
$ws.Cells.Item(17, 8).Value = 1796.875
$ws.Cells.Item(17, 10).Value = 1796.875
$ws.Cells.Item(17, 12).Value = 5390.625
$ws.Cells.Item(17, 14).Value = -5726.625
# Row 106 (diff hunk @@ -5886,22 +5886,22 @@)
$ws.Cells.Item(106, 8).Value = 5373524
$ws.Cells.Item(106, 9).Value = 7266534
$ws.Cells.Item(106, 11).Value = 7266534
$ws.Cells.Item(106, 13).Value = -7265903
# Row 125 (diff hunk @@ -6832,25 +6832,25 @@)
$ws.Cells.Item(125, 8).Value = 6459.8184
$ws.Cells.Item(125, 9).Value = 12887.667
$ws.Cells.Item(125, 10).Value = 4049.375
$ws.Cells.Item(125, 11).Value = 115989.003
$ws.Cells.Item(125, 12).Value = 36444.375
$ws.Cells.Item(125, 13).Value = -113529.003
$ws.Cells.Item(125, 14).Value = -41364.375
# Row 132 (diff hunk @@ -7178,22 +7178,22 @@)
$ws.Cells.Item(132, 8).Value = 4473.7812
$ws.Cells.Item(132, 9).Value = 4295.0356
$ws.Cells.Item(132, 11).Value = 12885.1068
$ws.Cells.Item(132, 13).Value = -10355.1068
# Row 137 (diff hunk @@ -7429,22 +7429,22 @@)
$ws.Cells.Item(137, 8).Value = 9525.23
$ws.Cells.Item(137, 9).Value = 11929.526
$ws.Cells.Item(137, 11).Value = 35788.578
$ws.Cells.Item(137, 13).Value = -33238.578
# Row 138 (diff hunk @@ -7481,22 +7481,22 @@)
$ws.Cells.Item(138, 8).Value = 2320.7466
$ws.Cells.Item(138, 9).Value = 961.55884
$ws.Cells.Item(138, 11).Value = 2884.67652
$ws.Cells.Item(138, 13).Value = 2255.32348

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2 (diff hunk @@ -7780,22 +7780,22 @@)
$ws.Cells.Item(2, 8).Value = 40693.48
$ws.Cells.Item(2, 9).Value = 556
$ws.Cells.Item(2, 11).Value = 556
$ws.Cells.Item(2, 13).Value = -443
# Row 32 (diff hunk @@ -9235,25 +9235,25 @@)
$ws.Cells.Item(32, 8).Value = 6644.1787
$ws.Cells.Item(32, 9).Value = 6668.037
$ws.Cells.Item(32, 10).Value = 6000
$ws.Cells.Item(32, 11).Value = 6668.037
$ws.Cells.Item(32, 12).Value = 6000
$ws.Cells.Item(32, 13).Value = -6381.037
$ws.Cells.Item(32, 14).Value = -6574
# Row 61 (diff hunk @@ -10641,25 +10641,25 @@)
$ws.Cells.Item(61, 8).Value = 3947.1428
$ws.Cells.Item(61, 10).Value = 4996.25
$ws.Cells.Item(61, 12).Value = 4996.25
$ws.Cells.Item(61, 14).Value = -5420.25
# Row 63 (diff hunk @@ -10739,22 +10739,22 @@)
$ws.Cells.Item(63, 8).Value = 1404.8572
$ws.Cells.Item(63, 9).Value = 1404.8572
$ws.Cells.Item(63, 11).Value = 1404.8572
$ws.Cells.Item(63, 13).Value = -718.8571999999999
# Row 66 (diff hunk @@ -10880,22 +10880,22 @@)
$ws.Cells.Item(66, 8).Value = 1404.8572
$ws.Cells.Item(66, 9).Value = 1404.8572
$ws.Cells.Item(66, 11).Value = 7024.286
$ws.Cells.Item(66, 13).Value = -3592.286
# Row 82 (diff hunk @@ -11646,22 +11646,19 @@)
$ws.Cells.Item(82, 8).Value = 0
$ws.Cells.Item(82, 10).Value = 0
$ws.Cells.Item(82, 12).Value = 0
$ws.Cells.Item(82, 14).Value = ""
# Row 85 (diff hunk @@ -11796,22 +11793,19 @@)
$ws.Cells.Item(85, 8).Value = 0
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 14).Value = ""
# Row 97 (diff hunk @@ -12384,25 +12378,25 @@)
$ws.Cells.Item(97, 8).Value = 7146962.5
$ws.Cells.Item(97, 9).Value = 5670
$ws.Cells.Item(97, 10).Value = 22223024
$ws.Cells.Item(97, 11).Value = 5670
$ws.Cells.Item(97, 12).Value = 22223024
$ws.Cells.Item(97, 13).Value = -5174
$ws.Cells.Item(97, 14).Value = -22224016
# Row 102 (diff hunk @@ -12632,25 +12626,25 @@)
$ws.Cells.Item(102, 8).Value = 10469.125
$ws.Cells.Item(102, 10).Value = 6635.067
$ws.Cells.Item(102, 12).Value = 6635.067
$ws.Cells.Item(102, 14).Value = -9879.066999999999
# Row 116 (diff hunk @@ -13306,22 +13300,22 @@)
$ws.Cells.Item(116, 8).Value = 40693.48
$ws.Cells.Item(116, 9).Value = 556
$ws.Cells.Item(116, 11).Value = 556
$ws.Cells.Item(116, 13).Value = 1738
# Row 132 (diff hunk @@ -14066,25 +14060,25 @@)
$ws.Cells.Item(132, 8).Value = 3229.1562
$ws.Cells.Item(132, 9).Value = 2642.7
$ws.Cells.Item(132, 10).Value = 4206.5835
$ws.Cells.Item(132, 11).Value = 7928.099999999999
$ws.Cells.Item(132, 12).Value = 12619.7505
$ws.Cells.Item(132, 13).Value = -5398.099999999999
$ws.Cells.Item(132, 14).Value = -17679.7505
# Row 136 (diff hunk @@ -14268,25 +14262,25 @@)
$ws.Cells.Item(136, 8).Value = 3947.1428
$ws.Cells.Item(136, 10).Value = 4996.25
$ws.Cells.Item(136, 12).Value = 14988.75
$ws.Cells.Item(136, 14).Value = -20088.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3 (diff hunk @@ -14702,22 +14696,22 @@)
$ws.Cells.Item(3, 8).Value = 40693.48
$ws.Cells.Item(3, 9).Value = 556
$ws.Cells.Item(3, 11).Value = 556
$ws.Cells.Item(3, 13).Value = -442
# Row 20 (diff hunk @@ -15529,25 +15523,25 @@)
$ws.Cells.Item(20, 8).Value = 3694.4211
$ws.Cells.Item(20, 9).Value = 2140.0908
$ws.Cells.Item(20, 10).Value = 5831.625
$ws.Cells.Item(20, 11).Value = 2140.0908
$ws.Cells.Item(20, 12).Value = 5831.625
$ws.Cells.Item(20, 13).Value = -1893.0908
$ws.Cells.Item(20, 14).Value = -6325.625
# Row 86 (diff hunk @@ -18745,22 +18739,22 @@)
$ws.Cells.Item(86, 8).Value = 6889.1177
$ws.Cells.Item(86, 9).Value = 10932.223
$ws.Cells.Item(86, 11).Value = 10932.223
$ws.Cells.Item(86, 13).Value = -9809.223
# Row 89 (diff hunk @@ -18898,22 +18892,22 @@)
$ws.Cells.Item(89, 8).Value = 6889.1177
$ws.Cells.Item(89, 9).Value = 10932.223
$ws.Cells.Item(89, 11).Value = 54661.115
$ws.Cells.Item(89, 13).Value = -49045.115
# Row 105 (diff hunk @@ -19688,22 +19682,22 @@)
$ws.Cells.Item(105, 8).Value = 130451.25
$ws.Cells.Item(105, 9).Value = 333870
$ws.Cells.Item(105, 11).Value = 333870
$ws.Cells.Item(105, 13).Value = -332123
# Row 107 (diff hunk @@ -19789,22 +19783,22 @@)
$ws.Cells.Item(107, 8).Value = 3797.6155
$ws.Cells.Item(107, 9).Value = 3737
$ws.Cells.Item(107, 11).Value = 3737
$ws.Cells.Item(107, 13).Value = -1817
# Row 123 (diff hunk @@ -20555,22 +20549,22 @@)
$ws.Cells.Item(123, 8).Value = 97000
$ws.Cells.Item(123, 10).Value = 97000
$ws.Cells.Item(123, 12).Value = 97000
$ws.Cells.Item(123, 14).Value = -106800
# Row 134 (diff hunk @@ -21070,22 +21064,22 @@)
$ws.Cells.Item(134, 8).Value = 8176.5
$ws.Cells.Item(134, 9).Value = 8176.5
$ws.Cells.Item(134, 11).Value = 24529.5
$ws.Cells.Item(134, 13).Value = -21994.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16 (diff hunk @@ -22233,25 +22227,25 @@)
$ws.Cells.Item(16, 8).Value = 1466.7727
$ws.Cells.Item(16, 9).Value = 1203.7858
$ws.Cells.Item(16, 10).Value = 1927
$ws.Cells.Item(16, 11).Value = 1203.7858
$ws.Cells.Item(16, 12).Value = 1927
$ws.Cells.Item(16, 13).Value = -916.7858000000001
$ws.Cells.Item(16, 14).Value = -2501
# Row 31 (diff hunk @@ -22953,25 +22947,25 @@)
$ws.Cells.Item(31, 8).Value = 9888.200000000001
$ws.Cells.Item(31, 9).Value = 12240.077
$ws.Cells.Item(31, 10).Value = 5520.4287
$ws.Cells.Item(31, 11).Value = 12240.077
$ws.Cells.Item(31, 12).Value = 5520.4287
$ws.Cells.Item(31, 13).Value = -11945.077
$ws.Cells.Item(31, 14).Value = -6110.4287
# Row 34 (diff hunk @@ -23100,25 +23094,25 @@)
$ws.Cells.Item(34, 8).Value = 9888.200000000001
$ws.Cells.Item(34, 9).Value = 12240.077
$ws.Cells.Item(34, 10).Value = 5520.4287
$ws.Cells.Item(34, 11).Value = 12240.077
$ws.Cells.Item(34, 12).Value = 5520.4287
$ws.Cells.Item(34, 13).Value = -12038.077
$ws.Cells.Item(34, 14).Value = -5924.4287
# Row 58 (diff hunk @@ -24255,25 +24249,25 @@)
$ws.Cells.Item(58, 8).Value = 2435.5862
$ws.Cells.Item(58, 9).Value = 2447.68
$ws.Cells.Item(58, 10).Value = 2360
$ws.Cells.Item(58, 11).Value = 2447.68
$ws.Cells.Item(58, 12).Value = 2360
$ws.Cells.Item(58, 13).Value = -2244.68
$ws.Cells.Item(58, 14).Value = -2766
# Row 99 (diff hunk @@ -26261,22 +26255,22 @@)
$ws.Cells.Item(99, 8).Value = 29517584
$ws.Cells.Item(99, 9).Value = 39354110
$ws.Cells.Item(99, 11).Value = 39354110
$ws.Cells.Item(99, 13).Value = -39352612
# Row 113 (diff hunk @@ -26950,25 +26944,25 @@)
$ws.Cells.Item(113, 8).Value = 1466.7727
$ws.Cells.Item(113, 9).Value = 1203.7858
$ws.Cells.Item(113, 10).Value = 1927
$ws.Cells.Item(113, 11).Value = 1203.7858
$ws.Cells.Item(113, 12).Value = 1927
$ws.Cells.Item(113, 13).Value = 966.2141999999999
$ws.Cells.Item(113, 14).Value = -6267
# Row 126 (diff hunk @@ -27578,22 +27572,22 @@)
$ws.Cells.Item(126, 8).Value = 29517584
$ws.Cells.Item(126, 9).Value = 39354110
$ws.Cells.Item(126, 11).Value = 118062330
$ws.Cells.Item(126, 13).Value = -118059860
# Row 132 (diff hunk @@ -27866,22 +27860,22 @@)
$ws.Cells.Item(132, 8).Value = 3108.3
$ws.Cells.Item(132, 9).Value = 3108.3
$ws.Cells.Item(132, 11).Value = 9324.900000000001
$ws.Cells.Item(132, 13).Value = -6794.900000000001
# Row 136 (diff hunk @@ -28065,25 +28059,25 @@)
$ws.Cells.Item(136, 8).Value = 2435.5862
$ws.Cells.Item(136, 9).Value = 2447.68
$ws.Cells.Item(136, 10).Value = 2360
$ws.Cells.Item(136, 11).Value = 7343.039999999999
$ws.Cells.Item(136, 12).Value = 7080
$ws.Cells.Item(136, 13).Value = -4793.039999999999
$ws.Cells.Item(136, 14).Value = -12180

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 24 (diff hunk @@ -36725,22 +36719,22 @@)
$ws.Cells.Item(24, 8).Value = 1148577.1
$ws.Cells.Item(24, 9).Value = 4000000
$ws.Cells.Item(24, 11).Value = 4000000
$ws.Cells.Item(24, 13).Value = -3999827
# Row 27 (diff hunk @@ -36875,22 +36869,19 @@)
$ws.Cells.Item(27, 8).Value = 0
$ws.Cells.Item(27, 10).Value = 0
$ws.Cells.Item(27, 12).Value = 0
$ws.Cells.Item(27, 14).Value = ""
# Row 70 (diff hunk @@ -38961,22 +38952,22 @@)
$ws.Cells.Item(70, 8).Value = 6952344
$ws.Cells.Item(70, 9).Value = 10107430
$ws.Cells.Item(70, 11).Value = 10107430
$ws.Cells.Item(70, 13).Value = -10107160
# Row 73 (diff hunk @@ -39105,22 +39096,22 @@)
$ws.Cells.Item(73, 8).Value = 6952344
$ws.Cells.Item(73, 9).Value = 10107430
$ws.Cells.Item(73, 11).Value = 10107430
$ws.Cells.Item(73, 13).Value = -10106494
# Row 113 (diff hunk @@ -41044,22 +41035,22 @@)
$ws.Cells.Item(113, 8).Value = 9104.5
$ws.Cells.Item(113, 9).Value = 11452
$ws.Cells.Item(113, 11).Value = 11452
$ws.Cells.Item(113, 13).Value = -9282
# Row 122 (diff hunk @@ -41482,25 +41473,25 @@)
$ws.Cells.Item(122, 8).Value = 8610.826999999999
$ws.Cells.Item(122, 9).Value = 5782.5
$ws.Cells.Item(122, 10).Value = 17499.857
$ws.Cells.Item(122, 11).Value = 17347.5
$ws.Cells.Item(122, 12).Value = 52499.571
$ws.Cells.Item(122, 13).Value = -14897.5
$ws.Cells.Item(122, 14).Value = -57399.571

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16 (diff hunk @@ -43230,22 +43221,22 @@)
$ws.Cells.Item(16, 8).Value = 7973.25
$ws.Cells.Item(16, 9).Value = 8716.625
$ws.Cells.Item(16, 11).Value = 8716.625
$ws.Cells.Item(16, 13).Value = -8546.625
# Row 69 (diff hunk @@ -45791,20 +45782,23 @@)
$ws.Cells.Item(69, 8).Value = 105081.5
$ws.Cells.Item(69, 9).Value = 100000
$ws.Cells.Item(69, 11).Value = 100000
$ws.Cells.Item(69, 13).Value = -99189
# Row 72 (diff hunk @@ -45938,20 +45932,23 @@)
$ws.Cells.Item(72, 8).Value = 105081.5
$ws.Cells.Item(72, 9).Value = 100000
$ws.Cells.Item(72, 11).Value = 300000
$ws.Cells.Item(72, 13).Value = -295944
# Row 93 (diff hunk @@ -46934,25 +46931,25 @@)
$ws.Cells.Item(93, 8).Value = 8854.733
$ws.Cells.Item(93, 9).Value = 9460.166999999999
$ws.Cells.Item(93, 10).Value = 6433
$ws.Cells.Item(93, 11).Value = 9460.166999999999
$ws.Cells.Item(93, 12).Value = 6433
$ws.Cells.Item(93, 13).Value = -8212.166999999999
$ws.Cells.Item(93, 14).Value = -8929
# Row 100 (diff hunk @@ -47271,25 +47268,25 @@)
$ws.Cells.Item(100, 8).Value = 6499.3335
$ws.Cells.Item(100, 9).Value = 2999.6667
$ws.Cells.Item(100, 10).Value = 9999
$ws.Cells.Item(100, 11).Value = 2999.6667
$ws.Cells.Item(100, 12).Value = 9999
$ws.Cells.Item(100, 13).Value = -2458.6667
$ws.Cells.Item(100, 14).Value = -11081
# Row 122 (diff hunk @@ -48322,25 +48319,25 @@)
$ws.Cells.Item(122, 8).Value = 5473.6924
$ws.Cells.Item(122, 9).Value = 6265.8
$ws.Cells.Item(122, 10).Value = 2833.3333
$ws.Cells.Item(122, 11).Value = 18797.4
$ws.Cells.Item(122, 12).Value = 8499.999899999999
$ws.Cells.Item(122, 13).Value = -16347.4
$ws.Cells.Item(122, 14).Value = -13399.9999
# Row 132 (diff hunk @@ -48809,22 +48806,22 @@)
$ws.Cells.Item(132, 8).Value = 576968.9
$ws.Cells.Item(132, 9).Value = 831149.8
$ws.Cells.Item(132, 11).Value = 2493449.4
$ws.Cells.Item(132, 13).Value = -2490919.4

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 113 (diff hunk @@ -54751,22 +54748,22 @@)
$ws.Cells.Item(113, 8).Value = 2178.9473
$ws.Cells.Item(113, 9).Value = 1036.6364
$ws.Cells.Item(113, 11).Value = 3109.9092
$ws.Cells.Item(113, 13).Value = -939.9092000000001
# Row 131 (diff hunk @@ -55627,19 +55624,22 @@)
$ws.Cells.Item(131, 8).Value = 50000
$ws.Cells.Item(131, 10).Value = 50000
$ws.Cells.Item(131, 12).Value = 50000
$ws.Cells.Item(131, 14).Value = -60080
# Row 132 (diff hunk @@ -55673,22 +55673,22 @@)
$ws.Cells.Item(132, 8).Value = 9645.449000000001
$ws.Cells.Item(132, 9).Value = 10102.889
$ws.Cells.Item(132, 11).Value = 30308.667
$ws.Cells.Item(132, 13).Value = -27778.667
